# Generate Report for Handoff
# Updates Status from "In Translation" to "Ready for handoff" and refreshes
# the handoff timestamps on all three sheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: E2 (zh-cn status), F2 (de-de status), G2 (Latest HO Xliff Generate Date)
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-05 20:46:20"

# zh-cn sheet: C2 (Status), H2 (Latest Handoff Datetime)
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-05 20:46:16"

# de-de sheet: C2 (Status), H2 (Latest Handoff Datetime)
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-05 20:46:20"

# The Status columns grow to fit the longer "Ready for handoff" text
# (mirrors Excel auto-fitting the column after the content changed).
$wsOverview.Columns.Item(5).ColumnWidth = 16.3889
$wsOverview.Columns.Item(6).ColumnWidth = 16.3889
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3889
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3889
